# Applies updated odds values to Sheet1 per diff (Jogos_da_Semana_FlashScore_2024-10-12.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "M2" = 1.05
    "N2" = 11
    "AM14" = 1250
    "G17" = 2.7
    "I17" = 2.4
    "L17" = 3.2
    "M17" = 1.07
    "N17" = 9
    "AC17" = 9
    "AG17" = 7.5
    "AH17" = 12
    "AL17" = 34
    "AM17" = 301
    "AX17" = 15
    "AZ17" = 51
    "BB17" = 201
    "G18" = 2.57
    "J18" = 3.15
    "K18" = 2.12
    "O18" = 1.35
    "P18" = 2.95
    "S18" = 1.4
    "T18" = 2.7
    "U18" = 1.83
    "V18" = 1.87
    "X18" = 12
    "Z18" = 27
    "AE18" = 15.5
    "AG18" = 7.8
    "AJ18" = 26
    "AL18" = 35
    "AM18" = 700
    "AN18" = 4.45
    "AO18" = 13.5
    "AT18" = 2.7
    "AW18" = 4.4
    "G20" = 1.7
    "I20" = 5.25
    "J20" = 2.4
    "N20" = 7.5
    "AD20" = 7
    "AN20" = 3.5
    "AR20" = 51
    "AZ20" = 126
    "AT22" = 2.62
    "AT24" = 2.62
    "H30" = 4.4
    "I30" = 5.1
    "J30" = 1.93
    "K30" = 2.6
    "L30" = 4.85
    "T30" = 3.7
    "X30" = 9.25
    "AG30" = 23
    "AH30" = 40
    "AI30" = 17
    "AJ30" = 100
    "AO30" = 6.9
    "AQ30" = 18
    "AT30" = 3.7
    "AU30" = 6.9
    "AX30" = 26
    "G36" = 2.92
    "K36" = 2.05
    "L36" = 2.95
    "AB36" = 29
    "AF36" = 55
    "AM36" = 400
    "AT36" = 2.52
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

